$wb = $excel.ActiveWorkbook

# Finished Week 13 logging: M.Sargent is no longer on the roster -
# remove his entire column from both the Rushing and Receiving sheets.
foreach ($ws in $wb.Worksheets) {
    $ws.Range("H1").EntireColumn.Delete()
}
